$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell D2 from number 1979 to text "invariant"
$ws.Range("D2").Value = "invariant"

# Update cell D3 from "XRND" to "invariant"
$ws.Range("D3").Value = "invariant"

# Update selection to D4 (to match sheetView selection change)
$ws.Range("D4").Select()
